{"js": "// Split the paragraph that reads \"joao\" into two runs (\"J\" + \"oao\"), and\n// add a new paragraph \"Giane\" right after it. The \"_GoBack\" bookmark that\n// used to sit on the \"joao\" paragraph now sits on the new \"Giane\" paragraph\n// (mirrors how Word leaves its last-edit bookmark on the most recently\n// touched paragraph).\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = body.paragraphs.items.find((p) => p.text === \"joao\");\nif (!target) {\n  throw new Error('Paragraph with text \"joao\" not found');\n}\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p w:rsidR=\"00556BD2\" w:rsidRDefault=\"00556BD2\"><w:r><w:t>J</w:t></w:r><w:r><w:t>oao</w:t></w:r></w:p>\n          <w:p><w:r><w:t>Giane</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.getRange(\"Whole\").insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Split the paragraph that reads \"joao\" into two runs (\"J\" + \"oao\"), and\n# add a new paragraph \"Giane\" right after it. The \"_GoBack\" bookmark that\n# used to sit on the \"joao\" paragraph now sits on the new \"Giane\" paragraph\n# (mirrors how Word leaves its last-edit bookmark on the most recently\n# touched paragraph).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"joao\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Paragraph with text 'joao' not found\"\n}\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p w:rsidR=\"00556BD2\" w:rsidRDefault=\"00556BD2\"><w:r><w:t>J</w:t></w:r><w:r><w:t>oao</w:t></w:r></w:p>\n          <w:p><w:r><w:t>Giane</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>'\n\n[void]$target.Range.InsertXML($xml)\n"}
